# Gantt chart update: refresh feature list + plan/actual numbers, clear
# trailing template rows, move selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project Planner")

# --- Row labels first (order matters for newly-introduced shared strings:
# they are appended to the shared-string table in the order they are first
# written, so rows 13/16/14/15/17 are set in that exact sequence to land on
# notes/class faq/share notes (markdown)/pomodoro timer/classroom invite link
# in the right slots). ---

$ws.Range("B5").Value = "Register"
$ws.Range("B6").Value = "Login/logout"
$ws.Range("B7").Value = "Account public/private setting on an options menu webpage"
$ws.Range("B8").Value = "Chat rooms for class"
$ws.Range("B9").Value = "Invite people to class/chat rooms"
$ws.Range("B10").Value = "Creating classrooms"
$ws.Range("B11").Value = "Delete Account (and change names of comments to deleted user)"
$ws.Range("B12").Value = "searchable accounts based on public and private setting"
$ws.Range("B13").Value = "notes"
$ws.Range("B16").Value = "class faq"
$ws.Range("B14").Value = "share notes (markdown)"
$ws.Range("B15").Value = " pomodoro timer"
$ws.Range("B17").Value = "classroom invite link"

# --- Row 5-17: plan/actual start + duration + percent-complete numbers ---

$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 6
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 4
$ws.Range("G5").Value = 1

$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 6
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 4
$ws.Range("G6").Value = 1

$ws.Range("C7").Value = 6
$ws.Range("D7").Value = 4
$ws.Range("E7").Value = 6
$ws.Range("F7").Value = 4
$ws.Range("G7").Value = 1

$ws.Range("C8").Value = 11
$ws.Range("D8").Value = 6
$ws.Range("E8").Value = 11
$ws.Range("F8").Value = 6
$ws.Range("G8").Value = 1

$ws.Range("C9").Value = 11
$ws.Range("D9").Value = 6
$ws.Range("E9").Value = 11
$ws.Range("F9").Value = 4
$ws.Range("G9").Value = 1

$ws.Range("C10").Value = 6
$ws.Range("D10").Value = 6
$ws.Range("E10").Value = 6
$ws.Range("F10").Value = 6
$ws.Range("G10").Value = 1

$ws.Range("C11").Value = 6
$ws.Range("D11").Value = 6
$ws.Range("E11").Value = 6
$ws.Range("F11").Value = 6
$ws.Range("G11").Value = 1

$ws.Range("C12").Value = 11
$ws.Range("D12").Value = 6
$ws.Range("E12").Value = 11
$ws.Range("F12").Value = 7
$ws.Range("G12").Value = 1

$ws.Range("C13").Value = 17
$ws.Range("D13").Value = 6
$ws.Range("E13").Value = 17
$ws.Range("F13").Value = 6
$ws.Range("G13").Value = 1

$ws.Range("C14").Value = 23
$ws.Range("D14").Value = 6
$ws.Range("E14").Value = 23
$ws.Range("F14").Value = 6
$ws.Range("G14").Value = 1

$ws.Range("C15").Value = 17
$ws.Range("D15").Value = 6
$ws.Range("E15").Value = 17
$ws.Range("F15").Value = 6
$ws.Range("G15").Value = 1

$ws.Range("C16").Value = 17
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 17
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 1

$ws.Range("C17").Value = 11
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 11
$ws.Range("F17").Value = 4
$ws.Range("G17").Value = 0.5

# --- Rows 18-30: clear the unused template rows ---
$ws.Range("B18:G30").ClearContents()

# --- Move the active selection ---
$ws.Range("AR11").Select()
